$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 659, pushing existing rows 659:671 down to 662:674.
$ws.Rows("659:661").Insert()

# Fill the 3 newly inserted rows (659-661) with the new weekly price entries.
# Columns A,B,C,E,F,G,H,I,J,K repeat the same market/product context used
# throughout this block; only D (date), L (quality), M-P (volume/prices),
# Q (unit), R (origin), S (price per kg) and T (kg/unit) vary per row.

$rows = @(
    @{ Row = 659; D = 44448; L = "1a amarillo"; M = 450; N = 4000;  O = 4000;  P = 4000;  Q = "`$/malla 14 kilos"; R = "Provincia de Quillota";  S = 286; T = 14 },
    @{ Row = 660; D = 44448; L = "1a amarillo"; M = 900; N = 3500;  O = 4000;  P = 3833;  Q = "`$/malla 14 kilos"; R = "Región de O'Higgins";    S = 274; T = 14 },
    @{ Row = 661; D = 44448; L = "2a amarillo"; M = 270; N = 3000;  O = 3000;  P = 3000;  Q = "`$/malla 14 kilos"; R = "Provincia de Quillota";  S = 214; T = 14 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
